$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B updates (activity counts)
$ws.Range("B2").Value = 1722.5
$ws.Range("B3").Value = 161
$ws.Range("B5").Value = 1163
$ws.Range("B6").Value = 789
$ws.Range("B7").Value = 284
$ws.Range("B9").Value = 1253
$ws.Range("B10").Value = 122
$ws.Range("B12").Value = 83

# Column C updates (hour) - all become 18
$ws.Range("C2").Value = 18
$ws.Range("C3").Value = 18
$ws.Range("C4").Value = 18
$ws.Range("C5").Value = 18
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 18
$ws.Range("C8").Value = 18
$ws.Range("C9").Value = 18
$ws.Range("C10").Value = 18
$ws.Range("C11").Value = 18
$ws.Range("C12").Value = 18
